# Add a new "Competitor" column (F) that names the opponent faced by
# Manchester United in each match of the season21_22 results table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header -------------------------------------------------------
$ws.Range("F1").Value = "Competitor"

# Copy the header formatting (bold, centered, bordered) from an
# existing header cell so F1 matches B1:E1 exactly.
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows ------------------------------------------------------
$competitors = @(
  "Leeds United",
  "Southampton",
  "Wolverhampton Wanderers",
  "Newcastle United",
  "Young Boys Berne",
  "West Ham United",
  "West Ham United",
  "Aston Villa",
  "Villarreal",
  "Everton",
  "Leicester City",
  "Atalanta",
  "Liverpool",
  "Tottenham Hotspur",
  "Atalanta",
  "Manchester City"
)

for ($i = 0; $i -lt $competitors.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = $competitors[$i]
}
